# Update the "想去人数" (want-to-go count) figures in column F that were
# refreshed by the latest gh-pages data generation run.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 315
$ws1.Range("F7").Value  = 5693
$ws1.Range("F9").Value  = 7681
$ws1.Range("F10").Value = 416
$ws1.Range("F13").Value = 3869
$ws1.Range("F25").Value = 5328
$ws1.Range("F27").Value = 2112
$ws1.Range("F30").Value = 7934
$ws1.Range("F33").Value = 2204
$ws1.Range("F42").Value = 1179
$ws1.Range("F46").Value = 2094
$ws1.Range("F47").Value = 131

# Sheet "本地生活"
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 577

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 577
$ws4.Range("F7").Value  = 315
$ws4.Range("F8").Value  = 5693
$ws4.Range("F9").Value  = 7681
$ws4.Range("F10").Value = 416
$ws4.Range("F11").Value = 3869
$ws4.Range("F25").Value = 5328
$ws4.Range("F27").Value = 2112
$ws4.Range("F30").Value = 7934
$ws4.Range("F33").Value = 2204
$ws4.Range("F40").Value = 1179
$ws4.Range("F44").Value = 2094
$ws4.Range("F45").Value = 131
